$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")
$ws.Activate()

# Delete row 2 (the transaction row with ID 34), shifting all rows below it up by one.
$ws.Rows.Item(2).Delete()

# Restore the selection to match what Excel leaves selected after a row delete
# (the new row 2, which took the place of the deleted row).
$ws.Range("A2:XFD2").Select()
